$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.393.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.952.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4765"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4024"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.55"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08402"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.053"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.967.37"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.551"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.136"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06580"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.47"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.808"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.410.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.288"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.178.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.82"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.904"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.147"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.36"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9758"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09595"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.444"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.660"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.586"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.893"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02319"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06212"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.245"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6184"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1914"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.347"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.23%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.04"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5939"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.053"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.388"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000328"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06802"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.27%  "
